# chore: update Sheets via scheduled runner
# Refresh market-price-derived columns (currentAveragePrice / NQ / HQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for the rows
# whose underlying item prices moved since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1771.0416
$ws.Range("J17").Value = 1809.762
$ws.Range("L17").Value = 5429.286
$ws.Range("N17").Value = -5765.286
$ws.Range("H43").Value = 74991.14
$ws.Range("J43").Value = 205155.8
$ws.Range("L43").Value = 205155.8
$ws.Range("N43").Value = -205293.8
$ws.Range("H69").Value = 4122.143
$ws.Range("I69").Value = 4337.3335
$ws.Range("J69").Value = 3960.75
$ws.Range("K69").Value = 13012.0005
$ws.Range("L69").Value = 11882.25
$ws.Range("M69").Value = -12138.0005
$ws.Range("N69").Value = -13630.25
$ws.Range("H72").Value = 4122.143
$ws.Range("I72").Value = 4337.3335
$ws.Range("J72").Value = 3960.75
$ws.Range("K72").Value = 39036.0015
$ws.Range("L72").Value = 35646.75
$ws.Range("M72").Value = -34668.0015
$ws.Range("N72").Value = -44382.75
$ws.Range("H86").Value = 7876.185
$ws.Range("I86").Value = 7585.95
$ws.Range("K86").Value = 7585.95
$ws.Range("M86").Value = -6462.95
$ws.Range("H89").Value = 7876.185
$ws.Range("I89").Value = 7585.95
$ws.Range("K89").Value = 37929.75
$ws.Range("M89").Value = -32313.75
$ws.Range("H93").Value = 52351
$ws.Range("J93").Value = 52351
$ws.Range("L93").Value = 52351
$ws.Range("N93").Value = -57343
$ws.Range("H100").Value = 3574.25
$ws.Range("I100").Value = 1200
$ws.Range("K100").Value = 1200
$ws.Range("M100").Value = -659
$ws.Range("H109").Value = 99922.5
$ws.Range("J109").Value = 99922.5
$ws.Range("L109").Value = 99922.5
$ws.Range("N109").Value = -102696.5
$ws.Range("H110").Value = 55021.5
$ws.Range("J110").Value = 55021.5
$ws.Range("L110").Value = 55021.5
$ws.Range("N110").Value = -63201.5
$ws.Range("H112").Value = 304602.94
$ws.Range("J112").Value = 386307.62
$ws.Range("L112").Value = 1158922.86
$ws.Range("N112").Value = -1161138.86
$ws.Range("H116").Value = 1115766.2
$ws.Range("I116").Value = 4485.087
$ws.Range("J116").Value = 4767118.5
$ws.Range("K116").Value = 4485.087
$ws.Range("L116").Value = 4767118.5
$ws.Range("M116").Value = -1043.087
$ws.Range("N116").Value = -4774002.5
$ws.Range("H138").Value = 13515275
$ws.Range("I138").Value = 1014.2083
$ws.Range("J138").Value = 20002120
$ws.Range("K138").Value = 3042.6249
$ws.Range("L138").Value = 60006360
$ws.Range("M138").Value = 2097.3751
$ws.Range("N138").Value = -60016640

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9452.09
$ws.Range("I32").Value = 6178.4443
$ws.Range("J32").Value = 23050.309
$ws.Range("K32").Value = 6178.4443
$ws.Range("L32").Value = 23050.309
$ws.Range("M32").Value = -5891.4443
$ws.Range("N32").Value = -23624.309
$ws.Range("H43").Value = 23348
$ws.Range("J43").Value = 23348
$ws.Range("L43").Value = 23348
$ws.Range("N43").Value = -23974
$ws.Range("H102").Value = 21489.732
$ws.Range("I102").Value = 9922.916999999999
$ws.Range("K102").Value = 9922.916999999999
$ws.Range("M102").Value = -8300.916999999999
$ws.Range("H104").Value = 58742.5
$ws.Range("J104").Value = 58742.5
$ws.Range("L104").Value = 58742.5
$ws.Range("N104").Value = -65730.5
$ws.Range("H110").Value = 1375.05
$ws.Range("I110").Value = 1338.9445
$ws.Range("K110").Value = 1338.9445
$ws.Range("M110").Value = 706.0554999999999
$ws.Range("H135").Value = 76383.71000000001
$ws.Range("J135").Value = 76383.71000000001
$ws.Range("L135").Value = 76383.71000000001
$ws.Range("N135").Value = -86523.71000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 64712
$ws.Range("J6").Value = 64712
$ws.Range("L6").Value = 64712
$ws.Range("N6").Value = -64938
$ws.Range("H51").Value = 77618.5
$ws.Range("J51").Value = 77618.5
$ws.Range("L51").Value = 77618.5
$ws.Range("N51").Value = -78600.5
$ws.Range("H53").Value = 38759.6
$ws.Range("J53").Value = 38772.25
$ws.Range("L53").Value = 38772.25
$ws.Range("N53").Value = -39920.25
$ws.Range("H54").Value = 7238
$ws.Range("I54").Value = 7238
$ws.Range("K54").Value = 7238
$ws.Range("M54").Value = -6754
$ws.Range("H86").Value = 1803.75
$ws.Range("J86").Value = 1965.5
$ws.Range("L86").Value = 1965.5
$ws.Range("N86").Value = -4211.5
$ws.Range("H89").Value = 1803.75
$ws.Range("J89").Value = 1965.5
$ws.Range("L89").Value = 9827.5
$ws.Range("N89").Value = -21059.5
$ws.Range("H99").Value = 1564739.9
$ws.Range("I99").Value = 1495.8889
$ws.Range("J99").Value = 2843757.8
$ws.Range("K99").Value = 1495.8889
$ws.Range("L99").Value = 2843757.8
$ws.Range("M99").Value = 2.111100000000079
$ws.Range("N99").Value = -2846753.8
$ws.Range("H105").Value = 2916.3333
$ws.Range("I105").Value = 2749.5
$ws.Range("K105").Value = 2749.5
$ws.Range("M105").Value = -1002.5
$ws.Range("H115").Value = 79060.60000000001
$ws.Range("J115").Value = 83325.5
$ws.Range("L115").Value = 83325.5
$ws.Range("N115").Value = -86459.5
$ws.Range("H119").Value = 72992
$ws.Range("J119").Value = 72992
$ws.Range("L119").Value = 72992
$ws.Range("N119").Value = -82668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 38735.637
$ws.Range("J9").Value = 38735.637
$ws.Range("L9").Value = 38735.637
$ws.Range("N9").Value = -39071.637
$ws.Range("H16").Value = 2116.7693
$ws.Range("J16").Value = 2999.6667
$ws.Range("L16").Value = 2999.6667
$ws.Range("N16").Value = -3573.6667
$ws.Range("H113").Value = 2116.7693
$ws.Range("J113").Value = 2999.6667
$ws.Range("L113").Value = 2999.6667
$ws.Range("N113").Value = -7339.6667
$ws.Range("H117").Value = 47794.668
$ws.Range("J117").Value = 47794.668
$ws.Range("L117").Value = 47794.668
$ws.Range("N117").Value = -56972.668
$ws.Range("H118").Value = 72057.42999999999
$ws.Range("J118").Value = 72057.42999999999
$ws.Range("L118").Value = 72057.42999999999
$ws.Range("N118").Value = -75371.42999999999
$ws.Range("H119").Value = 95592
$ws.Range("J119").Value = 95592
$ws.Range("L119").Value = 95592
$ws.Range("N119").Value = -105268
$ws.Range("H134").Value = 1774.6666
$ws.Range("I134").Value = 1195.28
$ws.Range("K134").Value = 3585.84
$ws.Range("M134").Value = -1050.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 3433.25
$ws.Range("J52").Value = 3433.25
$ws.Range("L52").Value = 10299.75
$ws.Range("N52").Value = -10831.75
$ws.Range("H54").Value = 4039
$ws.Range("I54").Value = 2600
$ws.Range("J54").Value = 4998.3335
$ws.Range("K54").Value = 7800
$ws.Range("L54").Value = 14995.0005
$ws.Range("M54").Value = -7241
$ws.Range("N54").Value = -16113.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 20694.166
$ws.Range("J58").Value = 18833
$ws.Range("L58").Value = 18833
$ws.Range("N58").Value = -19387
$ws.Range("H93").Value = 120000
$ws.Range("J93").Value = 120000
$ws.Range("L93").Value = 120000
$ws.Range("N93").Value = -123744
$ws.Range("H109").Value = 45929.8
$ws.Range("J109").Value = 45929.8
$ws.Range("L109").Value = 45929.8
$ws.Range("N109").Value = -48009.8
$ws.Range("H124").Value = 60000
$ws.Range("J124").Value = 60000
$ws.Range("L124").Value = 60000
$ws.Range("N124").Value = -69820
$ws.Range("H126").Value = 3060.5898
$ws.Range("I126").Value = 2467.2917
$ws.Range("J126").Value = 4009.8667
$ws.Range("K126").Value = 7401.875100000001
$ws.Range("L126").Value = 12029.6001
$ws.Range("M126").Value = -4931.875100000001
$ws.Range("N126").Value = -16969.6001
$ws.Range("H132").Value = 2060.3333
$ws.Range("I132").Value = 1663.0526
$ws.Range("J132").Value = 3570
$ws.Range("K132").Value = 4989.1578
$ws.Range("L132").Value = 10710
$ws.Range("M132").Value = -2459.1578
$ws.Range("N132").Value = -15770

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H118").Value = 94763.71000000001
$ws.Range("J118").Value = 94763.71000000001
$ws.Range("L118").Value = 94763.71000000001
$ws.Range("N118").Value = -98077.71000000001
$ws.Range("H136").Value = 6218.0625
$ws.Range("I136").Value = 7306.278
$ws.Range("J136").Value = 4818.9287
$ws.Range("K136").Value = 21918.834
$ws.Range("L136").Value = 14456.7861
$ws.Range("M136").Value = -19368.834
$ws.Range("N136").Value = -19556.7861
